# Applies the "tuned algorithm / added templates and generators1" edit to
# the "Tablib Dataset" sheet:
#   - fills in the previously-empty "description" column (C) with team/division
#     labels such as "1-E-MAJOR" for every existing row
#   - re-balances a handful of division boundaries (league/division columns D/E)
#     one row earlier than before
#   - appends three brand-new teams (rows 90-92) to the ITB division
#   - widens column C to fit the new labels and refreshes the frozen-pane/
#     selection view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cNames = @("1-E-MAJOR","2-E-MAJOR","3-E-MAJOR","4-E-MAJOR","5-E-MAJOR","6-E-MAJOR","7-E-MAJOR","8-E-MAJOR","1-W-MAJOR","2-W-MAJOR","3-W-MAJOR","4-W-MAJOR","5-W-MAJOR","1-E-MINOR","2-E-MINOR","3-E-MINOR","4-E-MINOR","5-E-MINOR","6-E-MINOR","7-E-MINOR","8-E-MINOR","9-E-MINOR","10-E-MINOR","1-W-MINOR","2-W-MINOR","3-W-MINOR","4-W-MINOR","5-W-MINOR","6-W-MINOR","1-E-PW","2-E-PW","3-E-PW","4-E-PW","5-E-PW","6-E-PW","7-E-PW","8-E-PW","9-E-PW","1-W-PW","2-W-PW","3-W-PW","4-W-PW","5-W-PW","6-W-PW","7-W-PW","1-E-CP","2-E-CP","3-E-CP","4-E-CP","5-E-CP","6-E-CP","7-E-CP","8-E-CP","9-E-CP","10-E-CP","11-E-CP","12-E-CP","1-W-CP","2-W-CP","3-W-CP","4-W-CP","5-W-CP","6-W-CP","7-W-CP","8-W-CP","9-W-CP","1-E-TBALL","2-E-TBALL","3-E-TBALL","4-E-TBALL","5-E-TBALL","6-E-TBALL","7-E-TBALL","8-E-TBALL","9-E-TBALL","10-E-TBALL","1-W-TBALL","2-W-TBALL","3-W-TBALL","4-W-TBALL","5-W-TBALL","6-W-TBALL","1-ITB","2-ITB","3-ITB","4-ITB","5-ITB","6-ITB","7-ITB","8-ITB","9-ITB")
$dVals  = @(1,1,1,1,1,1,1,1,1,1,1,1,1,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,2,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,5,6,6,6,6,6,6,6,6,6)
$eVals  = @(1,1,1,1,1,1,1,1,2,2,2,2,2,3,3,3,3,3,3,3,3,3,3,4,4,4,4,4,4,5,5,5,5,5,5,5,5,5,6,6,6,6,6,6,6,7,7,7,7,7,7,7,7,7,7,7,7,8,8,8,8,8,8,8,8,8,9,9,9,9,9,9,9,9,9,9,10,10,10,10,10,10,11,11,11,11,11,11,11,11,11)

$n = $cNames.Length   # 91 data rows: sheet rows 2..92

# --- Column C ("description"): new team/division labels, rows 2..92 -------
$cArr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) { $cArr[$i,0] = $cNames[$i] }
$ws.Range("C2:C92").Value2 = $cArr

# --- Columns D/E ("league"/"division"): updated boundaries, rows 2..92 ----
$dArr = New-Object 'object[,]' $n,1
$eArr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $dArr[$i,0] = $dVals[$i]
    $eArr[$i,0] = $eVals[$i]
}
$ws.Range("D2:D92").Value2 = $dArr
$ws.Range("E2:E92").Value2 = $eArr

# --- Columns A/B ("id"/"name") for the 3 brand-new rows (90..92) ----------
$newIds = New-Object 'object[,]' 3,1
$newIds[0,0] = 89
$newIds[1,0] = 90
$newIds[2,0] = 91
$ws.Range("A90:A92").Value2 = $newIds
$ws.Range("B90:B92").Value2 = $newIds

# --- Column width for the newly-populated column C -------------------------
$ws.Columns.Item(3).ColumnWidth = 12.77734375

# --- View state: scroll the frozen pane back to the top, then select A10 -
$ws.Range("A2").Select()
$ws.Range("A10").Select()

Write-Host "Applied description/division updates, 3 new rows, column width and view state."
